$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 64, shifting existing rows 64-70 down to 65-71.
$ws.Rows.Item(64).Insert()

# Fill the new row 64 with data. Columns A, B, C, E, F, G, H, I, R are the
# same constant values used throughout this block of rows; D, J, K, L, M,
# N, O, P, Q take the new values for this record.
$ws.Cells.Item(64, 1).Value = 4
$ws.Cells.Item(64, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(64, 3).Value = "Los Lagos"
$ws.Cells.Item(64, 4).Value = 44491
$ws.Cells.Item(64, 4).NumberFormat = $ws.Cells.Item(65, 4).NumberFormat
$ws.Cells.Item(64, 5).Value = 10
$ws.Cells.Item(64, 6).Value = 100112052
$ws.Cells.Item(64, 7).Value = "Albahaca"
$ws.Cells.Item(64, 8).Value = "Sin especificar"
$ws.Cells.Item(64, 9).Value = "Primera"
$ws.Cells.Item(64, 10).Value = 90
$ws.Cells.Item(64, 11).Value = 6000
$ws.Cells.Item(64, 12).Value = 6000
$ws.Cells.Item(64, 13).Value = 6000
$ws.Cells.Item(64, 14).Value = "$/paquete"
$ws.Cells.Item(64, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(64, 16).Value = 6000
$ws.Cells.Item(64, 17).Value = 1
$ws.Cells.Item(64, 18).Value = "Hortaliza"
